# "Generate Report for Archive"
# The localization run moved on from handoff into translation, so the
# Status/zh-cn/de-de cells that used to read "Ready for handoff" now read
# "In Translation". Because the new text is shorter, the columns that show
# it were re-autofit (narrower) when the report was regenerated.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$narrowWidth = 12.5   # produces the re-autofit column width after the shorter text

# --- Overview sheet: columns E (zh-cn) and F (de-de) carry the status text
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = $newStatus
$wsOverview.Range("E2:F3").EntireColumn.ColumnWidth = $narrowWidth

# --- zh-cn detail sheet: column C ("Status") carries the status text
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2:C3").Value = $newStatus
$wsZh.Range("C2:C3").EntireColumn.ColumnWidth = $narrowWidth

# --- de-de detail sheet: column C ("Status") carries the status text
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2:C3").Value = $newStatus
$wsDe.Range("C2:C3").EntireColumn.ColumnWidth = $narrowWidth
